$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "79.656.75"
$ws.Range("E2").Value = "  +4.12%  "
$ws.Range("D3").Value = "3.206.12"
$ws.Range("E3").Value = "  +5.05%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "205.50"
$ws.Range("E5").Value = "  +1.38%  "
$ws.Range("D6").Value = "635.39"
$ws.Range("E6").Value = "  +1.40%  "
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.01%  "
$ws.Range("D8").Value = "0.241"
$ws.Range("E8").Value = "  +15.08%  "
$ws.Range("E9").Value = "  +6.09%  "
$ws.Range("D10").Value = "3.204.19"
$ws.Range("E10").Value = "  +4.98%  "
$ws.Range("D11").Value = "0.584"
$ws.Range("E11").Value = "  +33.09%  "
$ws.Range("E12").Value = "  +3.11%  "
$ws.Range("D13").Value = "5.51"
$ws.Range("E13").Value = "  +7.52%  "
$ws.Range("B14").Value = "ShibaInu"
$ws.Range("C14").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D14").Value = "0.0000233"
$ws.Range("E14").Value = "  +20.15%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "3.789.93"
$ws.Range("E15").Value = "  +4.90%  "
$ws.Range("D16").Value = "32.02"
$ws.Range("E16").Value = "  +8.47%  "
$ws.Range("D17").Value = "79.522.40"
$ws.Range("E17").Value = "  +4.05%  "
$ws.Range("D18").Value = "3.204.11"
$ws.Range("E18").Value = "  +4.96%  "
$ws.Range("D19").Value = "14.59"
$ws.Range("E19").Value = "  +8.05%  "
$ws.Range("D20").Value = "2.99"
$ws.Range("E20").Value = "  +29.46%  "
$ws.Range("D21").Value = "9.20"
$ws.Range("E21").Value = "  +1.13%  "
$ws.Range("D22").Value = "430.33"
$ws.Range("E22").Value = "  +14.80%  "
$ws.Range("D23").Value = "5.13"
$ws.Range("E23").Value = "  +17.83%  "
$ws.Range("B24").Value = "Aptos"
$ws.Range("C24").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D24").Value = "11.29"
$ws.Range("E24").Value = "  +13.57%  "
$ws.Range("B25").Value = "WrappedeETH"
$ws.Range("C25").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D25").Value = "3.366.21"
$ws.Range("E25").Value = "  +5.04%  "
$ws.Range("D26").Value = "4.79"
$ws.Range("E26").Value = "  +8.17%  "
$ws.Range("D27").Value = "77.16"
$ws.Range("E27").Value = "  +4.85%  "
$ws.Range("E28").Value = "  +0.11%  "
$ws.Range("D29").Value = "0.0000120"
$ws.Range("E29").Value = "  +7.14%  "
$ws.Range("D30").Value = "9.07"
$ws.Range("E30").Value = "  +9.05%  "
$ws.Range("D31").Value = "0.999"
$ws.Range("E31").Value = "  +0.14%  "
$ws.Range("D32").Value = "1.50"
$ws.Range("E32").Value = "  +5.61%  "
$ws.Range("D33").Value = "525.04"
$ws.Range("E33").Value = "  +2.90%  "
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("D35").Value = "0.145"
$ws.Range("E35").Value = "  +28.47%  "
$ws.Range("D36").Value = "22.98"
$ws.Range("E36").Value = "  +9.86%  "
$ws.Range("D37").Value = "0.119"
$ws.Range("E37").Value = "  +11.63%  "
$ws.Range("E38").Value = "  +0.00%  "
$ws.Range("E39").Value = "  +5.16%  "
$ws.Range("D40").Value = "164.83"
$ws.Range("E40").Value = "  +1.07%  "
$ws.Range("D42").Value = "192.49"
$ws.Range("E42").Value = "  +1.09%  "
$ws.Range("E43").Value = "  +0.05%  "
$ws.Range("D44").Value = "5.55"
$ws.Range("E44").Value = "  +6.33%  "
$ws.Range("D45").Value = "0.830"
$ws.Range("E45").Value = "  -0.38%  "
$ws.Range("E46").Value = "  +8.03%  "
$ws.Range("D47").Value = "1.33"
$ws.Range("E47").Value = "  +2.00%  "
$ws.Range("D48").Value = "43.35"
$ws.Range("E48").Value = "  +2.68%  "
$ws.Range("D49").Value = "25.97"
$ws.Range("E49").Value = "  +15.48%  "
$ws.Range("E50").Value = "  +4.39%  "
$ws.Range("B51").Value = "dogwifhat"
$ws.Range("C51").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D51").Value = "2.53"
$ws.Range("E51").Value = "  +2.93%  "
